$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.775.80"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.859.35"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'244.67"
$ws.Range("D6").Value = "'0.6427"
$ws.Range("E6").Value = "  +4.36%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'48.03"
$ws.Range("E8").Value = "  +5.10%  "
$ws.Range("D9").Value = "'0.07540"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("D10").Value = "'0.2987"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").Value = "'24.59"
$ws.Range("E11").Value = "  +6.05%  "
$ws.Range("D12").Value = "'0.07690"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "1.871.07"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "'5.051"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "'0.6912"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "'83.98"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "'0.000009892"
$ws.Range("E17").Value = "  +10.43%  "
$ws.Range("D18").Value = "'6.142"
$ws.Range("E18").Value = "  +4.88%  "
$ws.Range("D19").Value = "29.801.53"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "2.111.44"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'236.69"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'12.66"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'158.61"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").Value = "'8.559"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'17.91"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'0.06259"
$ws.Range("E30").Value = "  +7.49%  "
$ws.Range("D31").Value = "'1.496"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'1.289"
$ws.Range("E32").Value = "  +5.57%  "
$ws.Range("D33").Value = "'4.153"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "'4.096"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'1.902"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("D37").Value = "'0.7310"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'2.822"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'0.01795"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "1.219.35"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "'6.305"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").Value = "'0.9209"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "2.015.37"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'102.11"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "'67.08"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").Value = "'0.00000000119"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'0.4071"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "'9.164"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'1.674"
$ws.Range("E51").Value = "  +5.83%  "
